{"js": "// Password Test Case SDET.docx edit:\n//  1) The last paragraph of the original \"2. Test Code (Python)\" section\n//     (currently split across 3 runs + 2 proofErr markers because of the\n//     red squiggly \"grammar\" marks Word inserted around\n//     print(f'Password: {pwd} - {message}')) is collapsed into one plain\n//     run of text.\n//  2) The entire duplicated \"2. NEW Test Code (Python)\" section that\n//     follows (heading + all of the \"address\" themed Python lines, down\n//     to the final empty paragraph) is removed completely, so the\n//     document ends right after the \"print(...)\" line above.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// Locate the two anchor paragraphs by their text instead of hard-coded\n// indices, so the script is resilient to minor paragraph-count drift.\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\nlet printIdx = -1;      // \"    print(f'Password: {pwd} - {message}')\"\nlet newHeadingIdx = -1; // \"2. NEW Test Code (Python)\"\nfor (let i = 0; i < items.length; i++) {\n  const t = items[i].text;\n  if (printIdx === -1 && t.indexOf(\"print(f'Password: {pwd} - {message}'\") !== -1) {\n    printIdx = i;\n  }\n  if (newHeadingIdx === -1 && t.indexOf(\"NEW\") !== -1 && t.indexOf(\"Test Code (Python)\") !== -1) {\n    newHeadingIdx = i;\n  }\n}\n\nif (printIdx === -1 || newHeadingIdx === -1) {\n  throw new Error(\"Could not locate anchor paragraphs (printIdx=\" + printIdx + \", newHeadingIdx=\" + newHeadingIdx + \")\");\n}\n\nconst lastIdx = items.length - 1; // final (always-present) paragraph mark of the body\n\n// Step 1: delete every paragraph from the duplicate heading up to (but\n// not including) the very last paragraph of the body -- the body must\n// always keep at least one paragraph, so that final mark can't be\n// removed directly. Deleting a paragraph's own Range (rather than\n// calling Paragraph.delete()) merges its mark forward into the next\n// paragraph without touching that next paragraph's text.\nfor (let i = lastIdx - 1; i >= newHeadingIdx; i--) {\n  items[i].getRange().delete();\n}\nawait context.sync();\n\n// Step 2: re-load -- the trailing empty paragraph (originally the very\n// last paragraph in the document) and the \"print(...)\" paragraph are now\n// adjacent. Give the trailing paragraph the same style as the\n// \"print(...)\" paragraph, then delete the mark between them so the\n// \"print(...)\" content becomes the new final paragraph of the body\n// (matching the target, which has no empty paragraph after it).\nconst refreshed = body.paragraphs;\nrefreshed.load(\"items/style\");\nawait context.sync();\n\nconst trailing = refreshed.items[refreshed.items.length - 1];\nconst printPara = refreshed.items[printIdx];\n\ntrailing.style = printPara.style;\n\nconst endOfPrint = printPara.getRange(Word.InsertLocation.end);\nconst startOfTrailing = trailing.getRange(Word.InsertLocation.start);\nendOfPrint.expandTo(startOfTrailing).delete();\nawait context.sync();\n\n// Step 3: the merged final paragraph still has the original run split\n// (an xml:space=\"preserve\" \"    \" run, then a \"grammar-start\" run with\n// the printed text, then a \"grammar-end\" run with the trailing \")\").\n// Office.js doesn't expose proofErr marks directly, but replacing the\n// whole paragraph range's text collapses it down to one run and drops\n// the now-orphaned proofErr elements.\nconst finalParagraphs = body.paragraphs;\nfinalParagraphs.load(\"items/text\");\nawait context.sync();\n\nconst finalPara = finalParagraphs.items[finalParagraphs.items.length - 1];\nfinalPara.getRange().insertText(\n  \"    print(f'Password: {pwd} - {message}')\",\n  Word.InsertLocation.replace\n);\nawait context.sync();\n", "ps1": "# Password Test Case SDET.docx edit:\n#  1) The last paragraph of the original \"2. Test Code (Python)\" section\n#     (currently split across 3 runs + 2 proofErr markers because of the\n#     red squiggly \"grammar\" marks Word inserted around\n#     print(f'Password: {pwd} - {message}')) is collapsed into one plain\n#     run of text.\n#  2) The entire duplicated \"2. NEW Test Code (Python)\" section that\n#     follows (heading + all of the \"address\" themed Python lines, down\n#     to the final empty paragraph) is removed completely, so the\n#     document ends right after the \"print(...)\" line above.\n\n$d = $word.ActiveDocument\n\n# Locate the two anchor paragraphs by their text instead of hard-coded\n# indices, so the script is resilient to minor paragraph-count drift.\n$printIdx = -1\n$newHeadingIdx = -1\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $t = $d.Paragraphs.Item($i).Range.Text\n    if ($printIdx -eq -1 -and $t -like \"*print(f'Password: {pwd} - {message}'*\") {\n        $printIdx = $i\n    }\n    if ($newHeadingIdx -eq -1 -and $t -like \"*NEW*\" -and $t -like \"*Test Code (Python)*\") {\n        $newHeadingIdx = $i\n    }\n}\n\nif ($printIdx -eq -1 -or $newHeadingIdx -eq -1) {\n    throw \"Could not locate anchor paragraphs (printIdx=$printIdx, newHeadingIdx=$newHeadingIdx)\"\n}\n\n# --- Step 1: merge the split runs of the \"print(...)\" paragraph into one run ---\n# (Clear first, then set -- setting identical text in one shot is treated\n# as a no-op and leaves the original run/proofErr split untouched.)\n$printPara = $d.Paragraphs.Item($printIdx)\n$printRange = $d.Range($printPara.Range.Start, $printPara.Range.End - 1)\n$mergedText = $printRange.Text\n$printRange.Text = \"\"\n$printPara2 = $d.Paragraphs.Item($printIdx)\n$printRange2 = $d.Range($printPara2.Range.Start, $printPara2.Range.End - 1)\n$printRange2.Text = $mergedText\n\n# --- Step 2: delete the whole duplicated \"2. NEW Test Code (Python)\" section ---\n# Span from the start of the duplicate heading through the end of the\n# very last paragraph's own text (its paragraph mark is excluded here\n# since the body must always retain a final paragraph mark).\n$lastIdx = $d.Paragraphs.Count\n$startPara = $d.Paragraphs.Item($newHeadingIdx)\n$endPara = $d.Paragraphs.Item($lastIdx - 1)\n$bigRange = $d.Range($startPara.Range.Start, $endPara.Range.End)\n$bigRange.Delete()\n\n# --- Step 3: remove the now-orphaned trailing empty paragraph so the ---\n# \"print(...)\" paragraph becomes the final paragraph of the document,\n# matching the target (no empty paragraph left after it).\n$printPara3 = $d.Paragraphs.Item($printIdx)\n$trailingPara = $d.Paragraphs.Item($d.Paragraphs.Count)\n$trailingPara.Range.Style = $printPara3.Range.Style\n$markRng = $d.Range($printPara3.Range.End - 1, $trailingPara.Range.Start)\n$markRng.Delete()\n"}
